$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'Address',
    'stars1zg8f6kshwssmtw2mnvd0a9wvlq78qrgv3g4a5q',
    'stars1ddkp0uc87x76zks73t8gxw2dz0r4wm2npz80sd',
    'stars1ej2es5fjztqjcd4pwa0zyvaevtjd2y5whwzd00',
    'stars1f4d3q86j2u8xl5h8n0arcp2k6uaegqu999vtwk',
    'stars1j0lq4ynh3ct6kr09uh8x3dau57dt6uzpw9f2m0',
    'stars157cey3pkggh8zyun3vmp4ak5wpllwr0scwny89',
    'stars12msfftw69xave94ytr4ngta74rnvhvhzc6jqv6',
    'stars1xmnsjghk24v62hdm6n92ynqljwuaelwx76taa3',
    'stars1ayp9ga949jfeu9d73gtge0mz29978eyr6um2es',
    'stars13lcdnj8dgmaug4k5azydxvrclvmvqxkq4sdkva',
    'stars1r3dupp4savup3pmmq6d5zfsc444stucy4t9kan',
    'stars1xm2pj03uglrzpcsypuzrl7hrlumzdy7g8zdhcy',
    'stars14am952nzggxntteuq5napl6e27kerx4hcaa0d8',
    'stars1km4ce0jeu6xqe0qd0pj22q9k3dzav7njsaln3v',
    'stars1wt29e7ec5f04xy6ez5m6agzrweh0vxsunjv7ce',
    'stars1ga27mkfs8weg429dpw4qcf3n6r8jzwk5xyf8fm',
    'stars1pa0392fvwr23dre5yqytvqwyezpknawjyv8yzv',
    'stars1vx8r6x3lh6vh8wd4z2gt42sm359087e5lrdw5q',
    'stars1trghmhd4d9qcf9d0lsyesxzj3aasqg9stcpx59',
    'stars13ul82cuykzfly3fps55cdsu56eejr86ndudx4s',
    'stars194v8uwee2fvs2s8fa5k7j03ktwc87h5yxz8m7c',
    'stars10flwce9rgdwmecen8v8glfcfll7mphyw8vum2e',
    'stars1ccfvra0ksy09d6wx559fnrd2msppgdce9mgar8',
    'stars1c84vfnsz7futvraup0mf2ascx7mvq2jvsswc6w',
    'stars1zhqrfu9w3sugwykef3rq8t0vlxkz72vw3a4d8a',
    'stars1w875kdjh3yxxdurrctu2ze7hjm8mdp9qjvaywp',
    'stars1r4lylz39yupmec3wk3ajur6k0w7ykkwvuwq62c',
    'stars1axw88ugkakle6g6ma7ev4f30lhzzyfaz3rus5u',
    'stars1vtpt9kwjywm3xpake4v8sjmh2uxqq57hpejak9',
    'stars1puxpleecjkpss83q79ua8ac8e5fwwt79azzm6k',
    'stars199ay07lscpghl3stkq2mtg40p6t4r2u6eqsl8w',
    'stars10yasutlku2ayktvhh2whrl7ht9hzlxwyaph20r',
    'stars13sj7q26k96n24ydx9wp59y7yl2av06daj329hl',
    'stars1jue58ngp2888mqj6h5pmw0an6a20p87warwwjj',
    'stars1c8y4ja3llw0cpz2amv5qf90mstntf6ljmzzsuk',
    'stars17kal5lvug02rnmx8fdukxe6l4eyk7j2ejff86j',
    'stars10ye6m8k9g39t327yuvnj86y2kq639epzp4h9jv',
    'stars10cznnk39d8g4l255jc5vg6zv98taedxsmsmunz',
    'stars15xa0cshe43u9zc85qc78ysl7myl05rt6n2l2uq',
    'stars14qxxnvute7k8uceqtcljurzl0sv3dt4gkwvygg',
    'stars12scnajjddcrkskht7mj7qplzxh7zkzfw9rtdsw',
    'stars105qeucljvgkgf4vqc0nsksaflg3y7n3dnp5wdv',
    'stars10tl8mxxdnz6laskppws3epuew6nhmq7jt6205s',
    'stars1yhqft6d2msmzpugdjtawsgdlwvgq3samp4qg0r',
    'stars1rxv6lnypkes9xd2k0ff2cgqqnaeu74zm3lm4ua',
    'stars1m7a7nva00p82xr0tssye052r8sxsxvcygzqxwh',
    'stars1jume25ttjlcaqqjzjjqx9humvze3vcc8qf2kwl',
    'stars1j3p24vjzvdprntyf082us5g2h0rlzmycaxzmml',
    'stars1q2ke5n8rgr8ncjyeeqed3xyhyneluxgesludek',
    'stars1furt7dtl3ss9qpa9c0t8e6wkh3q3pehg9nhpw2',
    'stars1djv30wag0mth94tw6e3qwltl6swhzcwfc04sgm',
    'stars1gd4xlwv34t0chqtcftnc49x2pc7atlcut8r689',
    'stars1th8ktg3ulh0pxl5ym7csmemsc00ntnzcg55myu',
    'stars1fvnn96kq7xex7t365ltva3msu6ewntxf6ur826',
    'stars1mm96gp2feu4mtv8mvm2sr62njqr3xxxae55mf5',
    'stars15308cc5wx6ea5l565h5ncuuyulwc73z2nmwed2',
    'stars1xwz6xj36hk6cvu2pj6zrcdu76hupya70f6ndyx',
    'stars1n4a7szqm4gu0xegntwy245fjq8tlfh4vppj0e6',
    'stars1y9rxa93ply4qdypk42nrsqg5j50wn9pmna708n',
    'stars139a4n6w6dhwv60dj2clgwm6r0q84gju28z9at0',
    'stars1g93ux6aj97ny7dh6g3u3eerv9nnvlh0dq428zc',
    'stars1dqc78lf5fu9e2krvj2zegdfvshu0l00k8rzynf',
    'stars1qhcrmfmgdt6e550a544lzup5htv0svnkgjjmyj',
    'stars1s8qx0zvz8yd6e4x0mqmqf7fr9vvfn622wtp3g3',
    'stars17y9ysn4uwusc0qv0d48wtc5rf4gnu6mqvjpvg9',
    'stars1v6fnklkj5av596xrv9sucdvcgk7llhra0nf24d',
    'stars1xrrfek3qpz4ja8ccrjklfy064asyqq9k8wpjl7',
    'stars1h8vvw7rxger9zc7e5gmytm9v66lft9xuu9whhr',
    'stars12c4xfww70fe4h37hncdea8lt80dny8aegvpt8c',
    'stars1zeqwmhk9hqypxcp4973au358avnmchdjknvezs',
    'stars1sm7c84n4fryl39490zf3gr9550cln6g7lm4cpz',
    'stars1cdzspjq58ggdfxp3cjm9ggmdynkxrsvzjllfzc',
    'stars167ucsyj55sxh3vfsa000qs27vddepl0y5qsrlc',
    'stars1e9ypzs3qgrkwzpstvw7z4ag96qzv9qtdlg8luf',
    'stars1kwstysea0p8stxe349hru6hfxc78zsjnd2wmdz',
    'stars1jp3htxl899qsxx9q5aqm7f80v0sdu0sfnrz9ew',
    'stars1djntlfl6fvk3unn9md7edqxhhpapha6z3vpqrl',
    'stars1sux3szs8gcgkqymrqwd82geztwmpg6680t7cjv',
    'stars1yu00g53fk729rz94n26rrd5z593674sgtatavl',
    'stars1prhnhdk2nl3xugdrm32jch9d8wrs9xgw7rf62v',
    'stars1purcdnexunusedgpj5q6ssqqunay34h5x2xlgy',
    'stars1xd04a8wf3ym4fvqz9nut4xu3ldt3stjx8wkt42',
    'stars1z926ax906k0ycsuckele6x5hh66e2m4mxyqnss',
    'stars1mnzdlcjlvua94l47h77j420zhsdyechx6ukrrf',
    'stars1v4qfy658q6h2ypgkztz46thdfj2fstv50p323g',
    'stars124w26nuf4hxe0u8ywk8rfxsl9vn6cgef9plkl0',
    'stars1yg4ztgrgnmn52dftk3xrxzxxt3s0e54erf08gn',
    'stars15heg4vvf24zs86zxr3gerpgu35ge3qjqggtr48',
    'stars1etwruy48q0j2s0vj3qpee6nu938fvde6zm3mn0',
    'stars12x9hund9zux7ynft6sd7jsqdnmphp6ghfpr8u7',
    'stars1tz3ujjx9cq77unk09jvd2mdxvgdd0kv2n4774c',
    'stars1kfcxl4lm3md25yj66qkp4lkntyd78jy4laq3jv',
    'stars1vuy3a8zqtewxqx2hy555k22pdu2yk0x72cwyga',
    'stars1njyvry0t3j5dy4rr6ar5zfglg3cy2e8uumq3nn',
    'stars1rpnqgfdr79zv7zn6d8y888u7vqvpgtr5cqdnuw',
    'stars1mhv799y0029rpjlu3j8vs8r23tyf7e49n4hmt3',
    'stars1s33zct2zhhaf60x4a90cpe9yquw99jj0qh8pdx',
    'stars1c6wu0lmqdk74k3qj4k4lyyhv8qyht0ajy42s4l',
    'stars199dm0ugv50k4qgcsxjc8p2rertstltxvkddcc6',
    'stars10wxn2lv29yqnw2uf4jf439kwy5ef00qdyvtgfh',
    'stars1258f0g6yea03a7hcnyfsysjjkdavy64ma8rlsy',
    'stars1mgtyjevwq5pg5l6u8mac7skngzynt2ppdmh0at',
    'stars1lrmt70vunnfwa2z6ykdnts2rz0kzqvz6yl6fnk'
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

Write-Output "Updated $($values.Count) rows"
